$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.486.57"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "1.648.72"
$ws.Range("E3").Value = "  +2.63%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  +0.42%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.005"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "304.65"
$ws.Range("E6").Value = "  +0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3771"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.52"
$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3654"
$ws.Range("E9").Value = "  +0.71%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.256"
$ws.Range("E10").Value = "  -1.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08146"
$ws.Range("E11").Value = "  -0.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.006"
$ws.Range("E12").Value = "  +0.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.02"
$ws.Range("E13").Value = "  +0.23%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.667"
$ws.Range("E14").Value = "  +0.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001267"
$ws.Range("E15").Value = "  +1.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.350"
$ws.Range("E16").Value = "  -0.25%  "

$ws.Range("D17").Value = "1.647.40"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.56"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06933"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.26"
$ws.Range("E20").Value = "  +0.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.576"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  -0.16%  "

$ws.Range("D23").Value = "23.497.50"
$ws.Range("E23").Value = "  +1.12%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.92"
$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.263"
$ws.Range("E25").Value = "  +5.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.440"
$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.32"
$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.86"
$ws.Range("E28").Value = "  +1.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.316"
$ws.Range("E29").Value = "  +0.78%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "136.63"
$ws.Range("E30").Value = "  +1.30%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.333"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.831.90"
$ws.Range("E32").Value = "  +3.12%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.960"
$ws.Range("E33").Value = "  +2.41%  "

$ws.Range("B34").Value = "FraxShare"
$ws.Range("C34").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.99"
$ws.Range("E34").Value = "  +6.10%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9672"
$ws.Range("E35").Value = "  +0.31%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02878"
$ws.Range("E36").Value = "  +3.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.328"
$ws.Range("E37").Value = "  +3.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2570"
$ws.Range("E38").Value = "  +1.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07300"
$ws.Range("E39").Value = "  -2.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08883"
$ws.Range("E40").Value = "  +0.93%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.381"
$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7183"
$ws.Range("E42").Value = "  +1.12%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.51"
$ws.Range("E43").Value = "  +3.92%  "

$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.68"
$ws.Range("E44").Value = "  +1.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6604"
$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.373"
$ws.Range("E46").Value = "  +1.65%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.50%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.019"
$ws.Range("E48").Value = "  +0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08006"
$ws.Range("E49").Value = "  +0.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.224"
$ws.Range("E50").Value = "  +1.46%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.52"
$ws.Range("E51").Value = "  -4.15%  "
